# Update NATMI ligand-receptor edge metrics with recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03397833333333333
$ws.Range("H2").Value = 0.101935
$ws.Range("I2").Value = 0.03987224921182536
$ws.Range("J2").Value = 0.03987224921182535
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.06449866666666666
$ws.Range("N2").Value = 0.193496
$ws.Range("O2").Value = 0.004525829983623641
$ws.Range("P2").Value = 0.004525829983623642
$ws.Range("Q2").Value = 0.002191557195555555
$ws.Range("R2").Value = 0.01972401476
$ws.Range("S2").Value = 0.0001804550209973933
$ws.Range("T2").Value = 0.0001804550209973933

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03397833333333333
$ws.Range("H3").Value = 0.101935
$ws.Range("I3").Value = 0.03987224921182536
$ws.Range("J3").Value = 0.03987224921182535
$ws.Range("O3").Value = 0.745188142173877
$ws.Range("P3").Value = 0.7451881421738772
$ws.Range("Q3").Value = 0.3608448485544444
$ws.Range("R3").Value = 3.24760363699
$ws.Range("S3").Value = 0.02971232731445397
$ws.Range("T3").Value = 0.02971232731445397

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03397833333333333
$ws.Range("H4").Value = 0.101935
$ws.Range("I4").Value = 0.03987224921182536
$ws.Range("J4").Value = 0.03987224921182535
$ws.Range("M4").Value = 3.566885000000001
$ws.Range("N4").Value = 10.700655
$ws.Range("O4").Value = 0.2502860278424993
$ws.Range("P4").Value = 0.2502860278424993
$ws.Range("Q4").Value = 0.1211968074916667
$ws.Range("R4").Value = 1.090771267425
$ws.Range("S4").Value = 0.009979466876373992
$ws.Range("T4").Value = 0.009979466876373992

# Row 5
$ws.Range("G5").Value = 0.7475459999999999
$ws.Range("I5").Value = 0.8772160811096247
$ws.Range("J5").Value = 0.8772160811096247
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.06449866666666666
$ws.Range("N5").Value = 0.193496
$ws.Range("O5").Value = 0.004525829983623641
$ws.Range("P5").Value = 0.004525829983623642
$ws.Range("Q5").Value = 0.048215720272
$ws.Range("R5").Value = 0.433941482448
$ws.Range("S5").Value = 0.003970130842002768
$ws.Range("T5").Value = 0.003970130842002768

# Row 6
$ws.Range("G6").Value = 0.7475459999999999
$ws.Range("I6").Value = 0.8772160811096247
$ws.Range("J6").Value = 0.8772160811096247
$ws.Range("O6").Value = 0.745188142173877
$ws.Range("P6").Value = 0.7451881421738772
$ws.Range("Q6").Value = 7.938827384828
$ws.Range("R6").Value = 71.44944646345199
$ws.Range("S6").Value = 0.6536910217671302
$ws.Range("T6").Value = 0.6536910217671303

# Row 7
$ws.Range("G7").Value = 0.7475459999999999
$ws.Range("I7").Value = 0.8772160811096247
$ws.Range("J7").Value = 0.8772160811096247
$ws.Range("M7").Value = 3.566885000000001
$ws.Range("N7").Value = 10.700655
$ws.Range("O7").Value = 0.2502860278424993
$ws.Range("P7").Value = 0.2502860278424993
$ws.Range("Q7").Value = 2.66641061421
$ws.Range("R7").Value = 23.99769552789
$ws.Range("S7").Value = 0.2195549285004916
$ws.Range("T7").Value = 0.2195549285004917

# Row 8
$ws.Range("G8").Value = 0.07065566666666666
$ws.Range("H8").Value = 0.211967
$ws.Range("I8").Value = 0.08291166967854992
$ws.Range("J8").Value = 0.0829116696785499
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.06449866666666666
$ws.Range("N8").Value = 0.193496
$ws.Range("O8").Value = 0.004525829983623641
$ws.Range("P8").Value = 0.004525829983623642
$ws.Range("Q8").Value = 0.004557196292444444
$ws.Range("R8").Value = 0.041014766632
$ws.Range("S8").Value = 0.0003752441206234803
$ws.Range("T8").Value = 0.0003752441206234803

# Row 9
$ws.Range("G9").Value = 0.07065566666666666
$ws.Range("H9").Value = 0.211967
$ws.Range("I9").Value = 0.08291166967854992
$ws.Range("J9").Value = 0.0829116696785499
$ws.Range("O9").Value = 0.745188142173877
$ws.Range("P9").Value = 0.7451881421738772
$ws.Range("Q9").Value = 0.7503526758575555
$ws.Range("R9").Value = 6.753174082717999
$ws.Range("S9").Value = 0.06178479309229278
$ws.Range("T9").Value = 0.06178479309229279

# Row 10
$ws.Range("G10").Value = 0.07065566666666666
$ws.Range("H10").Value = 0.211967
$ws.Range("I10").Value = 0.08291166967854992
$ws.Range("J10").Value = 0.0829116696785499
$ws.Range("M10").Value = 3.566885000000001
$ws.Range("N10").Value = 10.700655
$ws.Range("O10").Value = 0.2502860278424993
$ws.Range("P10").Value = 0.2502860278424993
$ws.Range("Q10").Value = 0.2520206375983333
$ws.Range("R10").Value = 2.268185738385
$ws.Range("S10").Value = 0.02075163246563365
$ws.Range("T10").Value = 0.02075163246563365

Write-Host "Updated NATMI edge metrics for rows 2-10 (columns E:T) with new TPM-based values."
